$d = $word.ActiveDocument

# Helper: run a Find/Replace scoped to a single paragraph's Range so that
# matches in other, textually-identical paragraphs are left untouched.
function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $paraIndex for text: $findText"
    }
}

# 1) First "Describe what happens..." prompt (bulleted list, numId=6).
#    Becomes "Now run the program..." and keeps a trailing space.
Replace-InParagraph 26 `
    "Describe what happens when you run the program. Was your prediction correct?" `
    "Now run the program. Was your prediction correct? If it wasn’t, explain what was wrong. "

# 2) "set the volume to 0 initially and then change it by 10?" question.
Replace-InParagraph 30 " and then " " then "
Replace-InParagraph 30 " it by 10?" " it by positive 10?"

# 3) "set the volume to 0 initially and then change it by 20?" question.
Replace-InParagraph 42 " the volume to 0 initially and then " " the volume to 0 initially then "
Replace-InParagraph 42 " it by 20?" " it by positive 20?"

# 4) Second "Describe what happens..." prompt (bulleted list, numId=7).
#    Becomes "Now run the program..." with NO trailing space this time.
Replace-InParagraph 62 `
    "Describe what happens when you run the program. Was your prediction correct?" `
    "Now run the program. Was your prediction correct? If it wasn’t, explain what was wrong."
